$wb = $excel.ActiveWorkbook
$capital = $wb.Worksheets("capital")
$holding = $wb.Worksheets("holding")

# Remove the obsolete DataDate column from the "capital" sheet, shifting
# everything else one column to the left.
$capital.Columns("A").Delete() | Out-Null

# Update the remembered selection on the "holding" sheet.
$holding.Select() | Out-Null
$holding.Range("F7").Select() | Out-Null

# Make "capital" the active sheet/tab and set its remembered selection.
$capital.Select() | Out-Null
$capital.Range("B5").Select() | Out-Null

Write-Output "done"
